# Update automàtic: dades i banners [2026-02-27 22:50]
# Refreshes DATA_EXTRACCIO timestamps and a handful of re-measured
# weather readings (humitat/temperatura/pressio/radiacio) on the
# Dades_Meteo sheet. Percentage readings ("45%", "90%", ...) are written
# with a leading apostrophe so Excel stores them as literal text (matching
# the source data, which is plain text, not a numeric percentage format).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-27 22:48:21'
$ws.Range("O2").Value = '5.3 °C'
$ws.Range("E3").Value = '2026-02-27 22:48:24'
$ws.Range("H3").Value = '''45%'
$ws.Range("N3").Value = '-0.2 °C 22:16 TU'
$ws.Range("O3").Value = '4.1 °C'
$ws.Range("E4").Value = '2026-02-27 22:48:26'
$ws.Range("E5").Value = '2026-02-27 22:48:29'
$ws.Range("H5").Value = '''47%'
$ws.Range("O5").Value = '4.5 °C'
$ws.Range("E6").Value = '2026-02-27 22:48:31'
$ws.Range("H6").Value = '''90%'
$ws.Range("E7").Value = '2026-02-27 22:48:33'
$ws.Range("J7").Value = '1024.6 hPa'
$ws.Range("E8").Value = '2026-02-27 22:48:36'
$ws.Range("E9").Value = '2026-02-27 22:48:38'
$ws.Range("E10").Value = '2026-02-27 22:48:39'
$ws.Range("O10").Value = '10.7 °C'
$ws.Range("E11").Value = '2026-02-27 22:48:40'
$ws.Range("E12").Value = '2026-02-27 22:48:41'
$ws.Range("E13").Value = '2026-02-27 22:48:42'
$ws.Range("O13").Value = '6.6 °C'
$ws.Range("E14").Value = '2026-02-27 22:48:43'
$ws.Range("E15").Value = '2026-02-27 22:48:46'
$ws.Range("O15").Value = '10.7 °C'
$ws.Range("E16").Value = '2026-02-27 22:48:48'
$ws.Range("K16").Value = '15.6 MJ/m2'
$ws.Range("O16").Value = '2.5 °C'
$ws.Range("E17").Value = '2026-02-27 22:48:51'
$ws.Range("N17").Value = '4.2 °C 22:24 TU'
$ws.Range("E18").Value = '2026-02-27 22:48:53'
$ws.Range("E19").Value = '2026-02-27 22:48:55'
$ws.Range("K19").Value = '14.5 MJ/m2'
$ws.Range("O19").Value = '10.2 °C'
$ws.Range("E20").Value = '2026-02-27 22:48:56'
$ws.Range("E21").Value = '2026-02-27 22:48:58'
$ws.Range("E22").Value = '2026-02-27 22:49:00'
$ws.Range("O22").Value = '1.4 °C'
$ws.Range("E23").Value = '2026-02-27 22:49:03'
$ws.Range("O23").Value = '3.5 °C'
$ws.Range("E24").Value = '2026-02-27 22:49:05'
$ws.Range("H24").Value = '''76%'
$ws.Range("O24").Value = '10.1 °C'
$ws.Range("E25").Value = '2026-02-27 22:49:08'
$ws.Range("H25").Value = '''35%'
$ws.Range("N25").Value = '1.9 °C 22:14 TU'
$ws.Range("O25").Value = '5.8 °C'
$ws.Range("E26").Value = '2026-02-27 22:49:10'
$ws.Range("E27").Value = '2026-02-27 22:49:13'
$ws.Range("O27").Value = '5.4 °C'
$ws.Range("E28").Value = '2026-02-27 22:49:15'
$ws.Range("E29").Value = '2026-02-27 22:49:18'
$ws.Range("E30").Value = '2026-02-27 22:49:20'
$ws.Range("E31").Value = '2026-02-27 22:49:22'
$ws.Range("E32").Value = '2026-02-27 22:49:25'
$ws.Range("E33").Value = '2026-02-27 22:49:27'
$ws.Range("J33").Value = '1023.5 hPa'
$ws.Range("E34").Value = '2026-02-27 22:49:30'
$ws.Range("E35").Value = '2026-02-27 22:49:32'
$ws.Range("E36").Value = '2026-02-27 22:49:35'
$ws.Range("O36").Value = '11.5 °C'
$ws.Range("E37").Value = '2026-02-27 22:49:37'
$ws.Range("O37").Value = '8.3 °C'
$ws.Range("E38").Value = '2026-02-27 22:49:39'
$ws.Range("E39").Value = '2026-02-27 22:49:42'
$ws.Range("N39").Value = '0.8 °C 22:15 TU'
$ws.Range("O39").Value = '4.4 °C'
$ws.Range("E40").Value = '2026-02-27 22:49:45'
$ws.Range("E41").Value = '2026-02-27 22:49:47'
$ws.Range("E42").Value = '2026-02-27 22:49:50'
$ws.Range("E43").Value = '2026-02-27 22:49:52'
$ws.Range("O43").Value = '9.2 °C'
$ws.Range("E44").Value = '2026-02-27 22:49:54'
$ws.Range("E45").Value = '2026-02-27 22:49:57'
$ws.Range("E46").Value = '2026-02-27 22:49:59'
$ws.Range("H46").Value = '''82%'
